$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''29.056.47'
$ws.Range('E2').Value = '''  -0.54%  '

$ws.Range('D3').Value = '''1.829.32'
$ws.Range('E3').Value = '''  -0.31%  '

$ws.Range('D4').Value = '''0.9989'
$ws.Range('E4').Value = '''  -0.18%  '

$ws.Range('D5').Value = '''241.46'
$ws.Range('E5').Value = '''  +0.43%  '

$ws.Range('D6').Value = '''0.6320'
$ws.Range('E6').Value = '''  -5.02%  '

$ws.Range('D7').Value = '''1.000'
$ws.Range('E7').Value = '''  -0.06%  '

$ws.Range('D8').Value = '''44.75'
$ws.Range('E8').Value = '''  +7.28%  '

$ws.Range('E9').Value = '''  +0.82%  '

$ws.Range('D10').Value = '''0.07338'
$ws.Range('E10').Value = '''  +0.08%  '

$ws.Range('E11').Value = '''  +1.25%  '

$ws.Range('D12').Value = '''0.07678'
$ws.Range('E12').Value = '''  -0.02%  '

$ws.Range('D13').Value = '''1.832.80'
$ws.Range('E13').Value = '''  +0.15%  '

$ws.Range('D14').Value = '''4.991'
$ws.Range('E14').Value = '''  +0.72%  '

$ws.Range('D15').Value = '''0.6632'
$ws.Range('E15').Value = '''  -0.18%  '

$ws.Range('D16').Value = '''82.04'
$ws.Range('E16').Value = '''  -1.12%  '

$ws.Range('D17').Value = '''6.065'
$ws.Range('E17').Value = '''  -0.35%  '

$ws.Range('D18').Value = '''0.000008654'
$ws.Range('E18').Value = '''  +4.91%  '

$ws.Range('D19').Value = '''28.988.71'
$ws.Range('E19').Value = '''  -0.71%  '

$ws.Range('D20').Value = '''2.079.97'
$ws.Range('E20').Value = '''  +0.07%  '

$ws.Range('D21').Value = '''12.41'
$ws.Range('E21').Value = '''  -0.19%  '

$ws.Range('D22').Value = '''224.15'
$ws.Range('E22').Value = '''  -0.44%  '

$ws.Range('E23').Value = '''  -0.11%  '

$ws.Range('D24').Value = '''7.133'
$ws.Range('E24').Value = '''  +0.30%  '

$ws.Range('D25').Value = '''1.000'
$ws.Range('E25').Value = '''  -0.17%  '

$ws.Range('E26').Value = '''  -1.54%  '

$ws.Range('D27').Value = '''8.455'
$ws.Range('E27').Value = '''  -1.87%  '

$ws.Range('D28').Value = '''0.1368'
$ws.Range('E28').Value = '''  -1.48%  '

$ws.Range('E29').Value = '''  -0.18%  '

$ws.Range('D30').Value = '''1.505'
$ws.Range('E30').Value = '''  -0.31%  '

$ws.Range('D31').Value = '''4.092'
$ws.Range('E31').Value = '''  -0.24%  '

$ws.Range('D32').Value = '''4.021'
$ws.Range('E32').Value = '''  -0.15%  '

$ws.Range('E33').Value = '''  +0.65%  '

$ws.Range('D34').Value = '''0.05301'
$ws.Range('E34').Value = '''  +0.50%  '

$ws.Range('D35').Value = '''0.7394'
$ws.Range('E35').Value = '''  -1.00%  '

$ws.Range('D36').Value = '''1.828'
$ws.Range('E36').Value = '''  -1.84%  '

$ws.Range('D37').Value = '''1.153'
$ws.Range('E37').Value = '''  +2.29%  '

$ws.Range('D38').Value = '''2.651'
$ws.Range('E38').Value = '''  -1.12%  '

$ws.Range('D39').Value = '''1.293.15'
$ws.Range('E39').Value = '''  -1.61%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = '''2.749'
$ws.Range('E40').Value = '''  +1.20%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '''0.01781'
$ws.Range('E41').Value = '''  -0.47%  '

$ws.Range('D42').Value = '''6.294'
$ws.Range('E42').Value = '''  +6.01%  '

$ws.Range('D43').Value = '''0.8935'
$ws.Range('E43').Value = '''  -2.60%  '

$ws.Range('E44').Value = '''  -0.23%  '

$ws.Range('D45').Value = '''102.60'
$ws.Range('E45').Value = '''  +0.87%  '

$ws.Range('D46').Value = '''1.977.25'

$ws.Range('D47').Value = '''0.00000000123'
$ws.Range('E47').Value = '''  -3.09%  '

$ws.Range('D48').Value = '''0.5138'
$ws.Range('E48').Value = '''  -0.54%  '

$ws.Range('D49').Value = '''64.12'
$ws.Range('E49').Value = '''  +0.93%  '

$ws.Range('D50').Value = '''1.728'
$ws.Range('E50').Value = '''  -1.80%  '

$ws.Range('D51').Value = '''0.05822'
$ws.Range('E51').Value = '''  -1.92%  '
